$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value2 = 26999.5
$ws.Range("J75").Value2 = 26999.5
$ws.Range("L75").Value2 = 26999.5
$ws.Range("N75").Value2 = -28871.5

$ws.Range("H78").Value2 = 26999.5
$ws.Range("J78").Value2 = 26999.5
$ws.Range("L78").Value2 = 80998.5
$ws.Range("N78").Value2 = -90358.5

$ws.Range("H82").Value2 = 5230.0625
$ws.Range("I82").Value2 = 1753.7778
$ws.Range("J82").Value2 = 9699.571
$ws.Range("K82").Value2 = 5261.3334
$ws.Range("L82").Value2 = 29098.713
$ws.Range("M82").Value2 = -4855.3334
$ws.Range("N82").Value2 = -29910.713

$ws.Range("H85").Value2 = 5230.0625
$ws.Range("I85").Value2 = 1753.7778
$ws.Range("J85").Value2 = 9699.571
$ws.Range("K85").Value2 = 5261.3334
$ws.Range("L85").Value2 = 29098.713
$ws.Range("M85").Value2 = -3857.3334
$ws.Range("N85").Value2 = -31906.713

$ws.Range("H114").Value2 = 39578
$ws.Range("J114").Value2 = 39578
$ws.Range("L114").Value2 = 39578
$ws.Range("N114").Value2 = -48256


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value2 = 19996
$ws.Range("J9").Value2 = 19996
$ws.Range("L9").Value2 = 19996
$ws.Range("N9").Value2 = -20336

$ws.Range("H20").Value2 = 19996
$ws.Range("J20").Value2 = 19996
$ws.Range("L20").Value2 = 19996
$ws.Range("N20").Value2 = -20536

$ws.Range("H109").Value2 = 25641.762
$ws.Range("J109").Value2 = 25641.762
$ws.Range("L109").Value2 = 25641.762
$ws.Range("N109").Value2 = -28415.762

$ws.Range("H122").Value2 = 7130.5713
$ws.Range("I122").Value2 = 0
$ws.Range("J122").Value2 = 7130.5713
$ws.Range("K122").Value2 = 0
$ws.Range("L122").Value2 = 21391.7139
$ws.Range("N122").Value2 = -26291.7139

$ws.Range("H132").Value2 = 2489.925
$ws.Range("I132").Value2 = 1308.6154
$ws.Range("J132").Value2 = 4683.7856
$ws.Range("K132").Value2 = 3925.8462
$ws.Range("L132").Value2 = 14051.3568
$ws.Range("M132").Value2 = -1395.8462
$ws.Range("N132").Value2 = -19111.3568


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value2 = 200.81818
$ws.Range("I64").Value2 = 172.4
$ws.Range("J64").Value2 = 224.5
$ws.Range("K64").Value2 = 172.4
$ws.Range("L64").Value2 = 224.5
$ws.Range("M64").Value2 = 52.59999999999999
$ws.Range("N64").Value2 = -674.5

$ws.Range("H67").Value2 = 200.81818
$ws.Range("I67").Value2 = 172.4
$ws.Range("J67").Value2 = 224.5
$ws.Range("K67").Value2 = 172.4
$ws.Range("L67").Value2 = 224.5
$ws.Range("M67").Value2 = 607.6
$ws.Range("N67").Value2 = -1784.5

$ws.Range("H134").Value2 = 3557.3667
$ws.Range("I134").Value2 = 1988.32
$ws.Range("J134").Value2 = 11402.6
$ws.Range("K134").Value2 = 5964.96
$ws.Range("L134").Value2 = 34207.8
$ws.Range("M134").Value2 = -3429.96
$ws.Range("N134").Value2 = -39277.8


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 490.51352
$ws.Range("I22").Value2 = 291.84616
$ws.Range("J22").Value2 = 960.0909
$ws.Range("K22").Value2 = 291.84616
$ws.Range("L22").Value2 = 960.0909
$ws.Range("M22").Value2 = 58.15384
$ws.Range("N22").Value2 = -1660.0909

$ws.Range("H81").Value2 = 27000
$ws.Range("J81").Value2 = 27000
$ws.Range("L81").Value2 = 27000
$ws.Range("N81").Value2 = -28996

$ws.Range("H84").Value2 = 27000
$ws.Range("J84").Value2 = 27000
$ws.Range("L84").Value2 = 81000
$ws.Range("N84").Value2 = -90984

$ws.Range("H117").Value2 = 0
$ws.Range("J117").Value2 = 0
$ws.Range("L117").Value2 = 0

$ws.Range("H132").Value2 = 3646.7576
$ws.Range("I132").Value2 = 3392.9048
$ws.Range("J132").Value2 = 4091
$ws.Range("K132").Value2 = 10178.7144
$ws.Range("L132").Value2 = 12273
$ws.Range("M132").Value2 = -7648.714399999999
$ws.Range("N132").Value2 = -17333


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value2 = 42.11111
$ws.Range("I8").Value2 = 42.11111
$ws.Range("K8").Value2 = 126.33333
$ws.Range("M8").Value2 = 12.66667000000001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value2 = 11190
$ws.Range("J5").Value2 = 11299
$ws.Range("L5").Value2 = 11299
$ws.Range("N5").Value2 = -11523

$ws.Range("H88").Value2 = 34466.668
$ws.Range("J88").Value2 = 34466.668
$ws.Range("L88").Value2 = 34466.668
$ws.Range("N88").Value2 = -35368.668

$ws.Range("H91").Value2 = 34466.668
$ws.Range("J91").Value2 = 34466.668
$ws.Range("L91").Value2 = 34466.668
$ws.Range("N91").Value2 = -37586.668

$ws.Range("H102").Value2 = 2313.9443
$ws.Range("I102").Value2 = 2000.8667
$ws.Range("J102").Value2 = 3879.3333
$ws.Range("K102").Value2 = 2000.8667
$ws.Range("L102").Value2 = 3879.3333
$ws.Range("M102").Value2 = -378.8667
$ws.Range("N102").Value2 = -7123.3333

$ws.Range("H122").Value2 = 4531.5454
$ws.Range("I122").Value2 = 3214.7144
$ws.Range("K122").Value2 = 9644.143199999999
$ws.Range("M122").Value2 = -7194.143199999999

$ws.Range("H126").Value2 = 4117.23
$ws.Range("I126").Value2 = 2973.913
$ws.Range("J126").Value2 = 5399.9756
$ws.Range("K126").Value2 = 8921.739
$ws.Range("L126").Value2 = 16199.9268
$ws.Range("M126").Value2 = -6451.739
$ws.Range("N126").Value2 = -21139.9268


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value2 = 37249.168
$ws.Range("J74").Value2 = 40699
$ws.Range("L74").Value2 = 40699
$ws.Range("N74").Value2 = -42695

$ws.Range("H77").Value2 = 37249.168
$ws.Range("J77").Value2 = 40699
$ws.Range("L77").Value2 = 122097
$ws.Range("N77").Value2 = -132081

$ws.Range("H122").Value2 = 6688
$ws.Range("I122").Value2 = 3701
$ws.Range("J122").Value2 = 9675
$ws.Range("K122").Value2 = 11103
$ws.Range("L122").Value2 = 29025
$ws.Range("M122").Value2 = -8653
$ws.Range("N122").Value2 = -33925


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value2 = 16460.223
$ws.Range("I29").Value2 = 3673.3333
$ws.Range("J29").Value2 = 22853.666
$ws.Range("K29").Value2 = 3673.3333
$ws.Range("L29").Value2 = 22853.666
$ws.Range("M29").Value2 = -3383.3333
$ws.Range("N29").Value2 = -23433.666

$ws.Range("H126").Value2 = 315136.84
$ws.Range("I126").Value2 = 1192.8096
$ws.Range("J126").Value2 = 822277.25
$ws.Range("K126").Value2 = 3578.4288
$ws.Range("L126").Value2 = 2466831.75
$ws.Range("M126").Value2 = -1108.4288
$ws.Range("N126").Value2 = -2471771.75


# Explicit cell removals (value cleared entirely, matching diff cell deletions)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M122").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N117").ClearContents()
